$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlTop = -4160

# Insert a new column before column A, shifting existing data (labels,
# RawActivations, PercActivations, totalActivation) one column to the right.
# This carries formatting along with the values (e.g. old B1:D1 -> C1:E1,
# old A2:A20 "index" labels -> B2:B20).
$ws.Range("A:A").Insert()

# New header for the inserted "segments" index column; match the bold,
# bordered, centered look of the other header cells (C1:E1).
$headerCell = $ws.Range("B1")
$headerCell.Value = "segments"
$headerCell.Font.Bold = $true
$headerCell.Borders.LineStyle = 1
$headerCell.HorizontalAlignment = $xlCenter
$headerCell.VerticalAlignment = $xlTop

# Fill the new column A (rows 2-20) with the 0-based segment index that
# corresponds to each label now sitting in column B, using the bold,
# bordered, centered "index column" look that column A used to have.
for ($i = 0; $i -lt 19; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $i
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlTop
}

# The labels (now in column B, rows 2-20) should look like plain data cells,
# not the bold/bordered/centered index style they inherited from old column A.
$ws.Range("B2:B20").Style = "Normal"
